# Updated cryptos list on Mon Feb 26 20:38:26 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the literal text even when it looks like a number (keeps
    # trailing zeros, e.g. '402.05' must stay text, not become 402.05),
    # then strip the quote-prefix style Excel applies for text-forced
    # numerics so the cell style is left untouched.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '54.329.69'
$ws.Range("E2").Value = '  +4.82%  '
$ws.Range("D3").Value = '3.176.68'
$ws.Range("E3").Value = '  +2.40%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") '402.05'
$ws.Range("E5").Value = '  +3.47%  '
Set-TextValue $ws.Range("D6") '109.07'
$ws.Range("E6").Value = '  +4.86%  '
Set-TextValue $ws.Range("D7") '0.550'
$ws.Range("E7").Value = '  +0.75%  '
Set-TextValue $ws.Range("D8") '0.999'
$ws.Range("E8").Value = '  -0.16%  '
Set-TextValue $ws.Range("D9") '0.620'
$ws.Range("E9").Value = '  +4.86%  '
Set-TextValue $ws.Range("D10") '38.87'
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("E11").Value = '  +1.49%  '
Set-TextValue $ws.Range("D12") '0.0882'
$ws.Range("E12").Value = '  +1.80%  '
$ws.Range("D13").Value = '3.671.06'
$ws.Range("E13").Value = '  +2.11%  '
Set-TextValue $ws.Range("D14") '19.21'
$ws.Range("E14").Value = '  +2.32%  '
Set-TextValue $ws.Range("D15") '8.05'
$ws.Range("E15").Value = '  +2.42%  '
Set-TextValue $ws.Range("D16") '1.07'
$ws.Range("E16").Value = '  +8.48%  '
$ws.Range("D17").Value = '3.181.73'
$ws.Range("E17").Value = '  +2.65%  '
Set-TextValue $ws.Range("D18") '10.50'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '54.270.45'
$ws.Range("E19").Value = '  +4.39%  '
$ws.Range("E20").Value = '  +3.66%  '
Set-TextValue $ws.Range("D21") '12.84'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").Value = '0.0₃0995'
$ws.Range("E22").Value = '  +2.39%  '
Set-TextValue $ws.Range("D23") '71.59'
$ws.Range("E23").Value = '  +1.66%  '
Set-TextValue $ws.Range("D24") '273.49'
$ws.Range("E24").Value = '  +1.59%  '
Set-TextValue $ws.Range("D25") '3.28'
$ws.Range("E25").Value = '  +4.16%  '
Set-TextValue $ws.Range("D26") '8.03'
$ws.Range("E26").Value = '  -2.28%  '
Set-TextValue $ws.Range("D27") '27.71'
$ws.Range("E27").Value = '  +2.35%  '
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  +3.09%  '
$ws.Range("E32").Value = '  +6.39%  '
$ws.Range("E33").Value = '  +9.74%  '
$ws.Range("E34").Value = '  +3.26%  '
Set-TextValue $ws.Range("D36") '50.59'
$ws.Range("E36").Value = '  +0.45%  '
Set-TextValue $ws.Range("D37") '3.66'
$ws.Range("E37").Value = '  +7.83%  '
Set-TextValue $ws.Range("D38") '0.999'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  +9.11%  '
$ws.Range("E40").Value = '  +10.34%  '
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("E42").Value = '  +1.90%  '
$ws.Range("E43").Value = '  +1.16%  '
Set-TextValue $ws.Range("D44") '130.36'
$ws.Range("E44").Value = '  +2.15%  '
Set-TextValue $ws.Range("D45") '0.118'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  +1.77%  '
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '2.087.83'
$ws.Range("E49").Value = '  +1.51%  '
Set-TextValue $ws.Range("D50") '0.0341'
$ws.Range("E50").Value = '  +6.25%  '
Set-TextValue $ws.Range("D51") '0.0505'
$ws.Range("E51").Value = '  +7.62%  '
